# Converts the five markdown-source headings in this guidance page into
# plain "FirstParagraph"/"BodyText" paragraphs prefixed with '#'/'##',
# matching the publish step that demotes Heading1/Heading2 styled
# paragraphs (and retargets the paragraph that used to follow them).
#
# Mapping (derived from the target OOXML):
#   Heading1 "Data Security and Privacy"                 -> FirstParagraph "#Data Security and Privacy"
#     following FirstParagraph paragraph                 -> BodyText
#   Heading2 "Why are security and privacy important?"   -> FirstParagraph "##Why are security and privacy important?"
#     following FirstParagraph paragraph                 -> BodyText
#   Heading2 "When this applies"                         -> BodyText "##When this applies"
#     following FirstParagraph paragraph                 -> BodyText
#   Heading2 "Contacts"                                  -> BodyText "##Contacts"
#     following FirstParagraph paragraph                 -> BodyText
#   Heading2 "Feedback"                                  -> FirstParagraph "##Feedback"
#     following FirstParagraph paragraph                 -> BodyText

$d = $word.ActiveDocument

function Get-ParagraphAtStart($start) {
    # $Range.Paragraphs.Item(...) objects returned from a Find-derived
    # range are unreliable in this host (empty .Range.Text) - walk the
    # document's own Paragraphs collection and match on start offset
    # instead, which is robust.
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($candidate.Range.Start -eq $start) {
            return $candidate
        }
    }
    return $null
}

function Convert-Heading {
    param(
        [string]$HeadingText,
        [string]$Prefix,
        [string]$HeadingNewStyle,
        [string]$NextNewStyle
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($HeadingText, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Heading not found: $HeadingText"
    }

    $startPos = $rng.Start

    # Prepend the markdown prefix to the heading's own text.
    $rng.Text = $Prefix + $rng.Text

    # The heading paragraph and the paragraph immediately following it
    # both get re-styled.
    $headingPara = Get-ParagraphAtStart $startPos
    $headingPara.Style = $HeadingNewStyle

    $nextPara = $headingPara.Next()
    $nextPara.Style = $NextNewStyle
}

Convert-Heading "Data Security and Privacy" "#" "First Paragraph" "Body Text"

Convert-Heading "Why are security and privacy important?" "##" "First Paragraph" "Body Text"

Convert-Heading "When this applies" "##" "Body Text" "Body Text"

Convert-Heading "Contacts" "##" "Body Text" "Body Text"

Convert-Heading "Feedback" "##" "First Paragraph" "Body Text"
